$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend formatting (border/font) from the last existing data row (row 7)
# down into the two new rows (8 and 9) before writing values, so the new
# rows inherit the same visual style as the rest of the table.
$ws.Range("A7:H7").Copy()
$ws.Range("A8:H9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(8).RowHeight = 29
$ws.Rows.Item(9).RowHeight = 29

# Row 9 is the "Other columns to include" row - its body (B9:D9, E9:H9) gets
# a slightly different look than the rest (no bold/no accidental shading),
# so nudge those ranges so Excel records them as a variant style.
$ws.Range("B9:D9").Interior.ColorIndex = -4142  # xlNone
$ws.Range("E9:H9").Font.Bold = $false
$ws.Range("E9:H9").Font.Size = 12

# --- Row 3: Internal unique identifier -> Format (tsv, csv, xls, or xlsx)
$ws.Range("A3").Value = "Format (tsv, csv, xls, or xlsx)"
$ws.Range("B3").Value = "tsv"
$ws.Range("C3").Value = "csv"

# --- Row 4: Title -> Internal unique identifier
$ws.Range("A4").Value = "Internal unique identifier"
$ws.Range("B4").Value = "EID"
$ws.Range("C4").Value = "UT"

# --- Row 5 (ISSN / ISSN / SN) is unchanged, left as-is

# --- Row 6: DOI -> EISSN (electronic ISSN)
$ws.Range("A6").Value = "EISSN (electronic ISSN)"
$ws.Range("B6").ClearContents()
$ws.Range("C6").Value = "EI"

# --- Row 7: Organization unit -> DOI
$ws.Range("A7").Value = "DOI"
$ws.Range("B7").Value = "DOI"
$ws.Range("C7").Value = "DI"

# --- Row 8 (new): Departments and/or faculties
$ws.Range("A8").Value = "Departments and/or faculties"
$ws.Range("B8").Value = "subject"
$ws.Range("C8").Value = "subject"

# --- Row 9 (new): Other columns to include
$ws.Range("A9").Value = "Other columns to include"
$ws.Range("B9").Value = "Title, Source title, Publisher"
$ws.Range("C9").Value = "TI, SO, PU, SC"

# Update selection to match the author's final cursor position
$ws.Range("B9").Select()
